$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> updates for column D (Price) and column E (Volume(1h)).
# Rows not listed for a given column are left unchanged.
$updates = @{
    2  = @{ D = "67.320.17";  E = "  +0.40%  " }
    3  = @{ D = "2.550.79";   E = "  -2.35%  " }
    4  = @{ D = "1.00";       E = "  -0.04%  " }
    5  = @{ D = "591.22";     E = "  +0.11%  " }
    6  = @{ D = "174.39";     E = "  +5.37%  " }
    7  = @{ E = "  +0.03%  " }
    8  = @{ E = "  -0.52%  " }
    9  = @{ D = "2.550.05";   E = "  -2.32%  " }
    10 = @{ D = "0.140";      E = "  +1.97%  " }
    11 = @{ E = "  +1.07%  " }
    12 = @{ E = "  -2.68%  " }
    13 = @{ E = "  -0.53%  " }
    14 = @{ D = "27.16";      E = "  -0.43%  " }
    15 = @{ E = "  -2.57%  " }
    16 = @{ E = "  -0.41%  " }
    17 = @{ D = "67.215.44";  E = "  +0.14%  " }
    18 = @{ D = "2.546.33";   E = "  -2.52%  " }
    19 = @{ D = "8.09";       E = "  +3.52%  " }
    20 = @{ D = "11.45";      E = "  -2.77%  " }
    21 = @{ D = "356.50";     E = "  +0.68%  " }
    22 = @{ D = "4.24";       E = "  -0.78%  " }
    23 = @{ D = "4.69";       E = "  +1.38%  " }
    24 = @{ D = "1.98";       E = "  +3.11%  " }
    25 = @{ E = "  -0.09%  " }
    26 = @{ D = "70.12";      E = "  +1.70%  " }
    27 = @{ D = "9.88";       E = "  -5.97%  " }
    28 = @{ D = "2.685.14";   E = "  -2.38%  " }
    29 = @{ E = "  +0.25%  " }
    30 = @{ E = "  +0.84%  " }
    31 = @{ D = "537.54";     E = "  -0.50%  " }
    32 = @{ E = "  +5.07%  " }
    33 = @{ E = "  +0.88%  " }
    34 = @{ E = "  -0.45%  " }
    35 = @{ E = "  -1.45%  " }
    36 = @{ D = "1.00";       E = "  -0.06%  " }
    37 = @{ D = "1.48";       E = "  -0.25%  " }
    38 = @{ D = "158.20";     E = "  +0.68%  " }
    39 = @{ E = "  -0.48%  " }
    40 = @{ E = "  +1.19%  " }
    41 = @{ E = "  -1.70%  " }
    42 = @{ D = "1.81";       E = "  +0.64%  " }
    43 = @{ E = "  +1.57%  " }
    44 = @{ D = "2.59";       E = "  +7.79%  " }
    45 = @{ E = "  -0.02%  " }
    46 = @{ D = "39.85";      E = "  -0.51%  " }
    47 = @{ D = "152.57";     E = "  +0.89%  " }
    48 = @{ D = "0.566";      E = "  -1.49%  " }
    49 = @{ D = "0.0₆0281";   E = "  -6.54%  " }
    50 = @{ E = "  -1.00%  " }
    51 = @{ D = "1.73";       E = "  +1.83%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        # Force text storage so numeric-looking price strings (e.g. "1.00",
        # "591.22") aren't silently coerced into numbers by Excel, matching
        # the source workbook's inline-string (text) cell type.
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $vals["E"]
    }
}
